# New weekly price-report row is inserted at row 219 of the data table,
# shifting all the existing rows (219..326) down by one (to 220..327).
# This mirrors the workbook's "latest entry goes on top" weekly update
# pattern (commit: "Fruta / hortaliza, semanal").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 219; everything below
# (old rows 219-326) shifts down to 220-327 and the sheet dimension
# grows from A1:R326 to A1:R327 automatically.
$ws.Rows("219:219").Insert()

# Populate the newly inserted row 219 with the new weekly record.
$ws.Range("A219").Value = 9
$ws.Range("B219").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C219").Value = "Metropolitana"
$ws.Range("D219").Value = 45089
$ws.Range("E219").Value = 13
$ws.Range("F219").Value = 100112003
$ws.Range("G219").Value = "Ajo"
$ws.Range("H219").Value = "Chino"
$ws.Range("I219").Value = "Primera"
$ws.Range("J219").Value = 520
$ws.Range("K219").Value = 15000
$ws.Range("L219").Value = 16000
$ws.Range("M219").Value = 15500
$ws.Range("N219").Value = "$/caja 10 kilos"
$ws.Range("O219").Value = "China"
$ws.Range("P219").Value = 1550
$ws.Range("Q219").Value = 10
$ws.Range("R219").Value = "Hortaliza"
